$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-4 so the unit/building/result values line up with
#     the new row layout (same text as before, values are unchanged in
#     substance, only the underlying shared-string slots shift around) ---
$ws.Range("B2").Value = "1单元"
$ws.Range("C2").Value = "2栋"
$ws.Range("E2").Value = "阳"

$ws.Range("B3").Value = "1单元"
$ws.Range("C3").Value = "2栋"
$ws.Range("E3").Value = "阳"

$ws.Range("B4").Value = "2单元"
$ws.Range("C4").Value = "1栋"
$ws.Range("E4").Value = "阳"

# --- Append two new positive-case rows ---
$ws.Range("A5").Value = "水草"
$ws.Range("B5").Value = "3单元"
$ws.Range("C5").Value = "2栋"
$ws.Range("D5").Value = 206
$ws.Range("E5").Value = "阳"
$ws.Range("F5").Value = 363636

$ws.Range("A6").Value = "王西宁"
$ws.Range("B6").Value = "4单元"
$ws.Range("C6").Value = "1栋"
$ws.Range("D6").Value = 101
$ws.Range("E6").Value = "阳"
$ws.Range("F6").Value = 123

# Copy the formatting used by the other "name" cells in column A (bold,
# centered, bordered header-style) onto the two new name cells.
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A6").PasteSpecial(-4122)
$excel.CutCopyMode = $false
